$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.669
$ws.Range("D10").Value = -8.066000000000001
$ws.Range("D12").Value = -6.725
$ws.Range("D18").Value = -8.319000000000001
$ws.Range("D25").Value = -7.675999999999999
$ws.Range("D37").Value = -8.166
$ws.Range("D55").Value = -8.218
$ws.Range("D68").Value = -7.228999999999999
$ws.Range("D77").Value = -8.150000000000002
$ws.Range("D78").Value = -8.171000000000001
$ws.Range("D79").Value = -8.068999999999999
$ws.Range("D80").Value = -7.989
$ws.Range("D81").Value = -7.754
$ws.Range("D82").Value = -8.17
$ws.Range("D84").Value = -8.442
$ws.Range("D101").Value = -7.991
$ws.Range("D102").Value = -8.298
